$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header cells P1 and Q1, matching the existing header style (bold, bordered, centered)
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15
$headerRange = $ws.Range("P1:Q1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$arr = New-Object 'object[,]' 24,16
$arr[0,0] = 0.257463097964532
$arr[0,1] = 0
$arr[0,2] = 0.0542464103497764
$arr[0,3] = 1.370482070515905
$arr[0,4] = 0.4154920926790453
$arr[0,5] = 0.3087818240933444
$arr[0,6] = 0.01300981890876352
$arr[0,7] = 0.002514754006794462
$arr[0,8] = 0.2402481991716172
$arr[0,9] = 0.2211169310015337
$arr[0,10] = 0
$arr[0,11] = 2.059257633578881
$arr[0,12] = 0
$arr[0,13] = 0
$arr[0,14] = 0.7769725187097123
$arr[0,15] = 1.108083529420298
$arr[1,0] = 0.2246422946948172
$arr[1,1] = 0
$arr[1,2] = 0.04772623219331962
$arr[1,3] = 1.197966391313443
$arr[1,4] = 0.3795704651027947
$arr[1,5] = 0.2815688872465572
$arr[1,6] = 0.009819827637526472
$arr[1,7] = 0.003111914995075793
$arr[1,8] = 0.231408937087572
$arr[1,9] = 0.2295725780981273
$arr[1,10] = 0
$arr[1,11] = 1.808619869206439
$arr[1,12] = 0
$arr[1,13] = 0
$arr[1,14] = 0.7849636033788201
$arr[1,15] = 1.031608714190483
$arr[2,0] = 0.2040570537338198
$arr[2,1] = 0
$arr[2,2] = 0.04371861471988581
$arr[2,3] = 1.092174280855033
$arr[2,4] = 0.3579506490399211
$arr[2,5] = 0.2652102302030244
$arr[2,6] = 0.008024639132514666
$arr[2,7] = 0.003562419889844293
$arr[2,8] = 0.2262223798667549
$arr[2,9] = 0.2349759201325325
$arr[2,10] = 0
$arr[2,11] = 1.655064555466225
$arr[2,12] = 0
$arr[2,13] = 0
$arr[2,14] = 0.7905566180360282
$arr[2,15] = 0.9858919506271633
$arr[3,0] = 0.1946335842596483
$arr[3,1] = 0
$arr[3,2] = 0.04210792929427498
$arr[3,3] = 1.049077603601745
$arr[3,4] = 0.3490508587573657
$arr[3,5] = 0.2584315559824546
$arr[3,6] = 0.00733048593731056
$arr[3,7] = 0.003843773184303068
$arr[3,8] = 0.224063624494562
$arr[3,9] = 0.2371742682415635
$arr[3,10] = 0
$arr[3,11] = 1.593976149966693
$arr[3,12] = 0
$arr[3,13] = 0
$arr[3,14] = 0.7933548980704117
$arr[3,15] = 0.9669272895185941
$arr[4,0] = 0.1919284626034283
$arr[4,1] = 0
$arr[4,2] = 0.04186954873117799
$arr[4,3] = 1.041918563126984
$arr[4,4] = 0.3473450945487073
$arr[4,5] = 0.2570743282931218
$arr[4,6] = 0.007217076193929362
$arr[4,7] = 0.003993542624560753
$arr[4,8] = 0.2235831241874067
$arr[4,9] = 0.2374767416804614
$arr[4,10] = 0
$arr[4,11] = 1.585573613655413
$arr[4,12] = 0
$arr[4,13] = 0
$arr[4,14] = 0.79426255776022
$arr[4,15] = 0.9630316031794166
$arr[5,0] = 0.20081963461638
$arr[5,1] = 0
$arr[5,2] = 0.04377676005751141
$arr[5,3] = 1.091583697745293
$arr[5,4] = 0.3571917630184203
$arr[5,5] = 0.2644733137427409
$arr[5,6] = 0.00801422351465364
$arr[5,7] = 0.003829090547006864
$arr[5,8] = 0.225851144190294
$arr[5,9] = 0.2348301002603808
$arr[5,10] = 0
$arr[5,11] = 1.65899827890695
$arr[5,12] = 0
$arr[5,13] = 0
$arr[5,14] = 0.7917936234215901
$arr[5,15] = 0.9835516697140463
$arr[6,0] = 0.242084236955364
$arr[6,1] = 0
$arr[6,2] = 0.05210554360232322
$arr[6,3] = 1.310949284089858
$arr[6,4] = 0.4021675757855263
$arr[6,5] = 0.2984694518538262
$arr[6,6] = 0.0118725835027263
$arr[6,7] = 0.003022158072648828
$arr[6,8] = 0.2366989035943803
$arr[6,9] = 0.2237835655132036
$arr[6,10] = 0
$arr[6,11] = 1.979110551777637
$arr[6,12] = 0
$arr[6,13] = 0
$arr[6,14] = 0.7812511105626854
$arr[6,15] = 1.078696696250177
$arr[7,0] = 0.32505905733764
$arr[7,1] = 0
$arr[7,2] = 0.06827520960280964
$arr[7,3] = 1.742979117355091
$arr[7,4] = 0.4953525978759643
$arr[7,5] = 0.3694145425761803
$arr[7,6] = 0.02089887043689531
$arr[7,7] = 0.001729019170193524
$arr[7,8] = 0.2606426267902293
$arr[7,9] = 0.2041139784359847
$arr[7,10] = 0
$arr[7,11] = 2.601706697224699
$arr[7,12] = 0
$arr[7,13] = 0
$arr[7,14] = 0.7639710230374419
$arr[7,15] = 1.279880431663628
$arr[8,0] = 0.3830353841491245
$arr[8,1] = 0
$arr[8,2] = 0.0811748120849316
$arr[8,3] = 1.960242400157426
$arr[8,4] = 0.5597584506818549
$arr[8,5] = 0.4168937782110049
$arr[8,6] = 0.02791974858045077
$arr[8,7] = 0.001386126329657067
$arr[8,8] = 0.2763143753045512
$arr[8,9] = 0.1907143159653577
$arr[8,10] = 0
$arr[8,11] = 3.048430592173332
$arr[8,12] = 0
$arr[8,13] = 0
$arr[8,14] = 0.759434804733921
$arr[8,15] = 1.413780911907878
$arr[9,0] = 0.4122735639353579
$arr[9,1] = 0
$arr[9,2] = 0.09601697146820243
$arr[9,3] = 1.269092207456993
$arr[9,4] = 0.5314327340949632
$arr[9,5] = 0.3806791504036227
$arr[9,6] = 0.04298614709841431
$arr[9,7] = 0.001806032129336721
$arr[9,8] = 0.255515955363137
$arr[9,9] = 0.1850354780034595
$arr[9,10] = 0
$arr[9,11] = 3.117634972433365
$arr[9,12] = 0
$arr[9,13] = 0
$arr[9,14] = 0.7945478843004707
$arr[9,15] = 1.293577415521042
$arr[10,0] = 0.4304065867277842
$arr[10,1] = 0
$arr[10,2] = 0.1062336575264879
$arr[10,3] = 0.772260580107428
$arr[10,4] = 0.4952574171301762
$arr[10,5] = 0.3423808720182819
$arr[10,6] = 0.07871157873649537
$arr[10,7] = 0.0017753205872717
$arr[10,8] = 0.2354775514269392
$arr[10,9] = 0.1834445775220006
$arr[10,10] = 0
$arr[10,11] = 3.073692180563171
$arr[10,12] = 0
$arr[10,13] = 0
$arr[10,14] = 0.8247922893306026
$arr[10,15] = 1.170452334576595
$arr[11,0] = 0.4359885643770838
$arr[11,1] = 0
$arr[11,2] = 0.113484688393811
$arr[11,3] = 0.395998123353273
$arr[11,4] = 0.4500498777536563
$arr[11,5] = 0.2993878321384358
$arr[11,6] = 0.1316413990087142
$arr[11,7] = 0.001771818275740245
$arr[11,8] = 0.2143535203355924
$arr[11,9] = 0.1845966922472977
$arr[11,10] = 0
$arr[11,11] = 2.94827562261807
$arr[11,12] = 0
$arr[11,13] = 0
$arr[11,14] = 0.8543645270133737
$arr[11,15] = 1.035029049673994
$arr[12,0] = 0.4331928777297662
$arr[12,1] = 0
$arr[12,2] = 0.1172459180672405
$arr[12,3] = 0.2055366135989658
$arr[12,4] = 0.4145402876316666
$arr[12,5] = 0.2674968339503323
$arr[12,6] = 0.1791769658262865
$arr[12,7] = 0.001866808945123033
$arr[12,8] = 0.1991736603320291
$arr[12,9] = 0.1866851726388168
$arr[12,10] = 0
$arr[12,11] = 2.82590964016731
$arr[12,12] = 0
$arr[12,13] = 0
$arr[12,14] = 0.8752350134084281
$arr[12,15] = 0.9355705210168423
$arr[13,0] = 0.4283829766728786
$arr[13,1] = 0
$arr[13,2] = 0.117479911968374
$arr[13,3] = 0.1682240055515543
$arr[13,4] = 0.4037756422816159
$arr[13,5] = 0.2584268390375115
$arr[13,6] = 0.1910926281606606
$arr[13,7] = 0.002017117133374491
$arr[13,8] = 0.1951253167701736
$arr[13,9] = 0.1878030611110431
$arr[13,10] = 0
$arr[13,11] = 2.780848211171559
$arr[13,12] = 0
$arr[13,13] = 0
$arr[13,14] = 0.8805952432856259
$arr[13,15] = 0.9078323047494337
$arr[14,0] = 0.4001689323618649
$arr[14,1] = 0
$arr[14,2] = 0.1101911224483985
$arr[14,3] = 0.1652027631474233
$arr[14,4] = 0.3855904847554683
$arr[14,5] = 0.2471839561745739
$arr[14,6] = 0.1761961752730627
$arr[14,7] = 0.002290032319959501
$arr[14,8] = 0.1926986061726623
$arr[14,9] = 0.1928465875838423
$arr[14,10] = 0
$arr[14,11] = 2.620508907689441
$arr[14,12] = 0
$arr[14,13] = 0
$arr[14,14] = 0.8749759180782064
$arr[14,15] = 0.8787420868883089
$arr[15,0] = 0.3799203571222165
$arr[15,1] = 0
$arr[15,2] = 0.1026277411229159
$arr[15,3] = 0.2504662988207187
$arr[15,4] = 0.390378378999003
$arr[15,5] = 0.2552650389717357
$arr[15,6] = 0.137642844795522
$arr[15,7] = 0.002478265938826496
$arr[15,8] = 0.1988500814788878
$arr[15,9] = 0.1956641074467393
$arr[15,10] = 0
$arr[15,11] = 2.562319995144463
$arr[15,12] = 0
$arr[15,13] = 0
$arr[15,14] = 0.8599661584552223
$arr[15,15] = 0.9086544743277898
$arr[16,0] = 0.3672378742383415
$arr[16,1] = 0
$arr[16,2] = 0.09408621952373863
$arr[16,3] = 0.4848615890570258
$arr[16,4] = 0.4162088589090231
$arr[16,5] = 0.282063524285995
$arr[16,6] = 0.08560355182264345
$arr[16,7] = 0.00230520631588238
$arr[16,8] = 0.2135189461695575
$arr[16,9] = 0.1968819861002489
$arr[16,10] = 0
$arr[16,11] = 2.581608687473107
$arr[16,12] = 0
$arr[16,13] = 0
$arr[16,14] = 0.8342839731870129
$arr[16,15] = 0.9961418755322313
$arr[17,0] = 0.3555565670351513
$arr[17,1] = 0
$arr[17,2] = 0.08582316355594344
$arr[17,3] = 0.9202872841222671
$arr[17,4] = 0.4569062121157259
$arr[17,5] = 0.3224729930970227
$arr[17,6] = 0.04315202952923158
$arr[17,7] = 0.002411729139096153
$arr[17,8] = 0.233935757064188
$arr[17,9] = 0.1970941417655943
$arr[17,10] = 0
$arr[17,11] = 2.672123414896305
$arr[17,12] = 0
$arr[17,13] = 0
$arr[17,14] = 0.8062192169830666
$arr[17,15] = 1.124580465973281
$arr[18,0] = 0.3581404203135605
$arr[18,1] = 0
$arr[18,2] = 0.07808042612076349
$arr[18,3] = 1.898799968616956
$arr[18,4] = 0.5403955054846676
$arr[18,5] = 0.4019900207814118
$arr[18,6] = 0.02594600260796431
$arr[18,7] = 0.002267319863605621
$arr[18,8] = 0.2709203627726708
$arr[18,9] = 0.1938785390190496
$arr[18,10] = 0
$arr[18,11] = 2.945479993350489
$arr[18,12] = 0
$arr[18,13] = 0
$arr[18,14] = 0.7648969946928048
$arr[18,15] = 1.370781081768115
$arr[19,0] = 0.4024460952712587
$arr[19,1] = 0
$arr[19,2] = 0.08643874456451073
$arr[19,3] = 2.224867814249606
$arr[19,4] = 0.6016186737511902
$arr[19,5] = 0.4501890046345238
$arr[19,6] = 0.03290286907933071
$arr[19,7] = 0.001974016835565529
$arr[19,8] = 0.2887168720959039
$arr[19,9] = 0.1841065867322556
$arr[19,10] = 0
$arr[19,11] = 3.311793156197155
$arr[19,12] = 0
$arr[19,13] = 0
$arr[19,14] = 0.7562385732754322
$arr[19,15] = 1.510559408455947
$arr[20,0] = 0.4354898644379546
$arr[20,1] = 0
$arr[20,2] = 0.09229594249290329
$arr[20,3] = 2.386687153391165
$arr[20,4] = 0.6397890593350866
$arr[20,5] = 0.4797043396183511
$arr[20,6] = 0.03735185884485714
$arr[20,7] = 0.001554881276029541
$arr[20,8] = 0.2994793119824806
$arr[20,9] = 0.1780335032192699
$arr[20,10] = 0
$arr[20,11] = 3.53749315050959
$arr[20,12] = 0
$arr[20,13] = 0
$arr[20,14] = 0.7513699609221547
$arr[20,15] = 1.59587305569994
$arr[21,0] = 0.4216340181247062
$arr[21,1] = 0
$arr[21,2] = 0.0890718475202732
$arr[21,3] = 2.30024923658155
$arr[21,4] = 0.6200909179549114
$arr[21,5] = 0.4646433627987534
$arr[21,6] = 0.03495138378759588
$arr[21,7] = 0.001455183170841146
$arr[21,8] = 0.2940939619173264
$arr[21,9] = 0.1813421590057676
$arr[21,10] = 0
$arr[21,11] = 3.41113947878523
$arr[21,12] = 0
$arr[21,13] = 0
$arr[21,14] = 0.7521907687346001
$arr[21,15] = 1.552557555949136
$arr[22,0] = 0.3623863674973222
$arr[22,1] = 0
$arr[22,2] = 0.07703056252653795
$arr[22,3] = 1.974497193023922
$arr[22,4] = 0.5458018928896422
$arr[22,5] = 0.4076853972634069
$arr[22,6] = 0.02642715880397262
$arr[22,7] = 0.001754487845055053
$arr[22,8] = 0.2737957451908386
$arr[22,9] = 0.1942017862731658
$arr[22,10] = 0
$arr[22,11] = 2.943599829352536
$arr[22,12] = 0
$arr[22,13] = 0
$arr[22,14] = 0.7595838184699844
$arr[22,15] = 1.388876829193009
$arr[23,0] = 0.2973137957647651
$arr[23,1] = 0
$arr[23,2] = 0.0640477563173647
$arr[23,3] = 1.625775334166661
$arr[23,4] = 0.4685264735709822
$arr[23,5] = 0.3486758786600461
$arr[23,6] = 0.01827567968906352
$arr[23,7] = 0.002454207865050684
$arr[23,8] = 0.2533125304225479
$arr[23,9] = 0.2090276457161759
$arr[23,10] = 0
$arr[23,11] = 2.441563591315628
$arr[23,12] = 0
$arr[23,13] = 0
$arr[23,14] = 0.7704291954441302
$arr[23,15] = 1.220398707944753

$target = $ws.Range("B2:Q25")
$target.Value2 = $arr
